$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Apply the 43 existing covid_deaths (column C) corrections ---
$cChanges = @(
    ,(96,10)
    ,(105,13)
    ,(110,14)
    ,(130,15)
    ,(140,21)
    ,(145,23)
    ,(159,12)
    ,(174,11)
    ,(216,11)
    ,(220,10)
    ,(305,6)
    ,(332,7)
    ,(336,8)
    ,(396,13)
    ,(412,6)
    ,(414,4)
    ,(424,4)
    ,(443,6)
    ,(1016,21)
    ,(1076,23)
    ,(1082,68)
    ,(1117,27)
    ,(1171,16)
    ,(1196,12)
    ,(1243,12)
    ,(1270,33)
    ,(1303,33)
    ,(1314,5)
    ,(1326,23)
    ,(1378,9)
    ,(1422,6)
    ,(1436,5)
    ,(1437,6)
    ,(1474,6)
    ,(1483,9)
    ,(1486,4)
    ,(1487,3)
    ,(1495,5)
    ,(1502,4)
    ,(1515,6)
    ,(1531,4)
    ,(1533,2)
    ,(1545,6)
)
foreach ($item in $cChanges) {
    $r = $item[0]
    $val = $item[1]
    $ws.Cells.Item($r, 3).Value = $val
}

# --- 2) Rewrite rows 1547-1624 (tail of the date range) with revised/added records ---
$newRows = @(
    ,(1547,44274,"40-49",1)
    ,(1548,44274,"50-59",1)
    ,(1549,44274,"60-69",6)
    ,(1550,44275,"60-69",1)
    ,(1551,44275,"80+",1)
    ,(1552,44276,"60-69",3)
    ,(1553,44276,"70-79",3)
    ,(1554,44276,"80+",4)
    ,(1555,44277,"50-59",1)
    ,(1556,44277,"60-69",2)
    ,(1557,44277,"70-79",2)
    ,(1558,44277,"80+",2)
    ,(1559,44278,"40-49",1)
    ,(1560,44278,"60-69",1)
    ,(1561,44278,"70-79",2)
    ,(1562,44278,"80+",1)
    ,(1563,44279,"60-69",1)
    ,(1564,44279,"70-79",3)
    ,(1565,44279,"80+",2)
    ,(1566,44280,"60-69",1)
    ,(1567,44280,"70-79",6)
    ,(1568,44280,"80+",2)
    ,(1569,44281,"30-39",1)
    ,(1570,44281,"50-59",1)
    ,(1571,44281,"60-69",3)
    ,(1572,44281,"70-79",1)
    ,(1573,44281,"80+",2)
    ,(1574,44282,"60-69",1)
    ,(1575,44282,"70-79",1)
    ,(1576,44282,"80+",5)
    ,(1577,44283,"40-49",1)
    ,(1578,44283,"50-59",1)
    ,(1579,44283,"60-69",1)
    ,(1580,44283,"70-79",3)
    ,(1581,44283,"80+",1)
    ,(1582,44284,"50-59",1)
    ,(1583,44284,"60-69",1)
    ,(1584,44284,"70-79",1)
    ,(1585,44284,"80+",1)
    ,(1586,44285,"60-69",5)
    ,(1587,44285,"80+",1)
    ,(1588,44286,"40-49",1)
    ,(1589,44286,"50-59",2)
    ,(1590,44286,"60-69",1)
    ,(1591,44286,"70-79",2)
    ,(1592,44286,"80+",2)
    ,(1593,44287,"50-59",2)
    ,(1594,44287,"60-69",1)
    ,(1595,44287,"70-79",1)
    ,(1596,44287,"80+",3)
    ,(1597,44288,"40-49",1)
    ,(1598,44288,"60-69",1)
    ,(1599,44288,"70-79",3)
    ,(1600,44288,"80+",1)
    ,(1601,44289,"60-69",1)
    ,(1602,44289,"70-79",3)
    ,(1603,44289,"80+",1)
    ,(1604,44290,"30-39",1)
    ,(1605,44290,"70-79",1)
    ,(1606,44290,"80+",2)
    ,(1607,44291,"80+",2)
    ,(1608,44292,"50-59",2)
    ,(1609,44292,"60-69",2)
    ,(1610,44293,"60-69",2)
    ,(1611,44293,"70-79",1)
    ,(1612,44293,"80+",3)
    ,(1613,44294,"50-59",1)
    ,(1614,44294,"70-79",3)
    ,(1615,44294,"80+",1)
    ,(1616,44295,"50-59",1)
    ,(1617,44295,"60-69",3)
    ,(1618,44295,"80+",1)
    ,(1619,44296,"60-69",1)
    ,(1620,44296,"70-79",3)
    ,(1621,44296,"80+",3)
    ,(1622,44297,"60-69",1)
    ,(1623,44297,"80+",3)
    ,(1624,44298,"60-69",1)
)
foreach ($item in $newRows) {
    $r = $item[0]
    $dateVal = $item[1]
    $ageGrp = $item[2]
    $count = $item[3]
    $ws.Cells.Item($r, 1).Value = $dateVal
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $ageGrp
    $ws.Cells.Item($r, 3).Value = $count
}